$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 17 with the new "Unigram LM with Jelinek-Mercer smoothing" results
$ws.Range("A17").Value = "Unigram LM with Jelinek-Mercer smoothing"
$ws.Range("B17").Value = 0.2301
$ws.Range("C17").Value = 0.344
$ws.Range("D17").Value = 0.2787

# Update the selected cell to A17 (matches the saved selection in the sheet view)
$ws.Range("A17").Select() | Out-Null
